$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 17: TestCaseID was cleared (A17 no longer holds a value)
$ws.Range("A17").ClearContents()

# Row 18: fill in the "get subscription details" result columns
$ws.Range("F18").Value = "Case Pass, the deducted amount is as per the calculation"
$ws.Range("G18").Value = 150
$ws.Range("H18").Value = "2024-08-28 21:13:20"

# Row 19: fill in the "get subscription details" result columns
$ws.Range("F19").Value = "Case Pass, the deducted amount is as per the calculation"
$ws.Range("G19").Value = 75
$ws.Range("H19").Value = "2024-08-28 21:54:08"

# Row 20: new test case added, plus result columns
$ws.Range("A20").Value = 19
$ws.Range("F20").Value = "Case Pass, the deducted amount is as per the calculation"
$ws.Range("G20").Value = 165
$ws.Range("H20").Value = "2024-08-28 22:15:06"
